$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.457.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.61%  '

$ws.Range("D3").Value = "'1.823.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.93%  '

$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.98%  '

$ws.Range("D5").Value = "'332.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.46%  '

$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.91%  '

$ws.Range("D7").Value = "'0.4567"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.14%  '

$ws.Range("D8").Value = "'0.3829"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.05%  '

$ws.Range("D9").Value = "'46.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.48%  '

$ws.Range("D10").Value = "'0.07849"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.27%  '

$ws.Range("D11").Value = "'0.9571"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.16%  '

$ws.Range("D12").Value = "'21.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.82%  '

$ws.Range("D13").Value = "'5.834"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.82%  '

$ws.Range("D14").Value = "'1.812.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.85%  '

$ws.Range("D15").Value = "'7.047"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.21%  '

$ws.Range("D16").Value = "'1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.01%  '

$ws.Range("D17").Value = "'89.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.36%  '

$ws.Range("D18").Value = "'0.06585"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.13%  '

$ws.Range("D19").Value = "'0.00001019"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.13%  '

$ws.Range("D20").Value = "'17.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.59%  '

$ws.Range("D21").Value = "'1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.03%  '

$ws.Range("D22").Value = "'27.436.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.70%  '

$ws.Range("D23").Value = "'5.278"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.97%  '

$ws.Range("D24").Value = "'10.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.35%  '

$ws.Range("D25").Value = "'2.262"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.78%  '

$ws.Range("D26").Value = "'158.73"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").Value = "'2.026.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.67%  '

$ws.Range("D28").Value = "'19.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.70%  '

$ws.Range("D29").Value = "'2.041"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.13%  '

$ws.Range("D30").Value = "'5.267"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.61%  '

$ws.Range("D31").Value = "'117.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.38%  '

$ws.Range("D32").Value = "'0.09336"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.63%  '

$ws.Range("D33").Value = "'0.9279"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.89%  '

$ws.Range("D34").Value = "'3.568"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.59%  '

$ws.Range("D35").Value = "'5.217"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.60%  '

$ws.Range("D36").Value = "'1.315"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.01%  '

$ws.Range("D37").Value = "'0.05899"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.81%  '

$ws.Range("D38").Value = "'0.02184"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.00%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = "'8.070"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.90%  '

$ws.Range("B40").Value = 'Frax'
$ws.Range("C40").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D40").Value = "'1.002"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.97%  '

$ws.Range("D41").Value = "'1.139"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.42%  '

$ws.Range("D42").Value = "'0.5719"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.66%  '

$ws.Range("D43").Value = "'0.1816"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.69%  '

$ws.Range("D44").Value = "'9.891"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.34%  '

$ws.Range("D45").Value = "'1.261"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.12%  '

$ws.Range("D46").Value = "'0.5384"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.76%  '

$ws.Range("D47").Value = "'11.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.47%  '

$ws.Range("D48").Value = "'1.879"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.85%  '

$ws.Range("D49").Value = "'0.06959"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.69%  '

$ws.Range("D50").Value = "'110.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.31%  '

$ws.Range("D51").Value = "'1.002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -33.17%  '
